$d = $word.ActiveDocument

# Track revisions while editing so that only the text runs we actually
# touch get rewritten; surrounding runs (with identical formatting) are
# left completely untouched instead of being coalesced into the edit.
$d.TrackRevisions = $true

# --- Change 1 -------------------------------------------------------
# Split "[Inventions et découvertes] Les piles atomiques" into two runs:
#   "[Science et technique] " + "Les piles atomiques"
$rng1 = $d.Content
$rng1.Find.Execute("[Inventions et découvertes] Les piles atomiques")
$splitPos = $rng1.Start + 28   # length of "[Inventions et découvertes] "... first run incl. trailing space
$firstRun = $d.Range($rng1.Start, $splitPos)
$firstRun.Text = "[Science et technique] "

# --- Change 2 -------------------------------------------------------
# Merge the separate " " run and the "uranium est introduit..." run
# into a single run " uranium est introduit sous forme de barres, ce
# qui permet facilement de l".
$rng2 = $d.Content
$rng2.Find.Execute(" uranium est introduit sous forme de barres, ce qui permet facilement de l")
$tmp = $d.Range($rng2.Start, $rng2.End)
$tmp.Text = "TEMPMARKER uranium est introduit sous forme de barres, ce qui permet facilement de l"
$markerRng = $d.Range($tmp.Start, $tmp.Start + 10)
$markerRng.Text = ""

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
